# Sprint 8 workbook update:
#  - Rename sheet "Sprint" -> "Sprint 8"
#  - Insert a new "HORAS TRABALHADAS" column into the table (Tabela2), between
#    "QUEM REALIZOU" and "STATUS"
#  - Fill in hour values for the new column (all "1 Hora" except the
#    "Criação da planilha de testes" row, which gets "1 Hora e 30 Minutos")
#  - Cosmetic: hide gridlines and update the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Remember the existing STATUS column (column D) contents before we touch
# anything, so we can move them one column to the right afterwards.
$statusVals = @()
for ($r = 1; $r -le 8; $r++) {
    $statusVals += $ws.Cells.Item($r, 4).Value2
}

# Grow the table by one column. The engine always appends the new column at
# the end of the table (column E), expanding the table range to A1:E8 and
# wiring up the associated header/data conditional-format (dxf) entries.
$col = $lo.ListColumns.Add()

# Approximate the original column widths (C/D match, E matches old D, F
# matches old E).
$ws.Columns.Item(4).ColumnWidth = 28.6
$ws.Columns.Item(5).ColumnWidth = 21.1
$ws.Columns.Item(6).ColumnWidth = 15.5

# Put the new "HORAS TRABALHADAS" column into D (pushing STATUS to E).
$ws.Range("D1").Value = "HORAS TRABALHADAS"

# Type the one-off value first, then the repeated one, matching the order
# the values were first introduced into the sheet.
$ws.Range("D5").Value = "1 Hora e 30 Minutos"
$ws.Range("D2").Value = "1 Hora"
$ws.Range("D3").Value = "1 Hora"
$ws.Range("D4").Value = "1 Hora"
$ws.Range("D6").Value = "1 Hora"
$ws.Range("D7").Value = "1 Hora"
$ws.Range("D8").Value = "1 Hora"

# Restore the STATUS header + values into their new home, column E.
for ($r = 1; $r -le 8; $r++) {
    $ws.Cells.Item($r, 5).Value = $statusVals[$r - 1]
}

# Rename the worksheet tab.
$ws.Name = "Sprint 8"

# Cosmetic view changes: hide gridlines and set the active selection.
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("D6:D9").Select() | Out-Null
